# "Add files via upload" — re-uploaded matrix now for the DESIGN phase
# instead of the REQUIREMENTS phase. Only the title cell (A1, merged
# A1:D2) actually changes text; everything else on the sheet is
# untouched. The selection/active cell in the saved file is also moved
# onto the title's merged range (A1:D2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet title (merged cell A1:D2) to reflect the new phase.
$ws.Range("A1").Value = "MATRIZ DE RIESGOS FASE DE DISEÑO"

# Leave the selection on the title range, matching the saved workbook.
$ws.Range("A1:D2").Select()
